$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (values + matching header style/format copied from H1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Data rows 2-19, column I and J values
$data = @(
    @{Row=2;  I=1; J=6},
    @{Row=3;  I=1; J=6},
    @{Row=4;  I=1; J=6},
    @{Row=5;  I=1; J=6},
    @{Row=6;  I=1; J=5},
    @{Row=7;  I=1; J=3},
    @{Row=8;  I=1; J=4},
    @{Row=9;  I=1; J=6},
    @{Row=10; I=1; J=6},
    @{Row=11; I=1; J=7},
    @{Row=12; I=1; J=5},
    @{Row=13; I=1; J=5},
    @{Row=14; I=1; J=4},
    @{Row=15; I=3; J=6},
    @{Row=16; I=8; J=8},
    @{Row=17; I=6; J=7},
    @{Row=18; I=5; J=6},
    @{Row=19; I=2; J=3}
)

foreach ($item in $data) {
    $ws.Cells.Item($item.Row, 9).Value = $item.I
    $ws.Cells.Item($item.Row, 10).Value = $item.J
}
